$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet, row 4 (53789b32-... entry)
$wsZhCn.Range("D4").Value = "2016-01-20 03:20:09"
$wsZhCn.Range("G4").Value = "2016-01-20 03:21:02"

# de-de sheet, row 4 (53789b32-... entry)
$wsDeDe.Range("D4").Value = "2016-01-20 03:20:19"
$wsDeDe.Range("G4").Value = "2016-01-20 03:21:20"
